{"js": "// The Orion activity guide repeats a campaign-dates sentence several times\n// throughout the document (once per language/section copy). The author\n// reworded it from:\n//   \"Waktu Kampanye rasi bintang Orion 2022: ...\"\n// to:\n//   \"Waktu Kampanye 2022 untuk rasi bintang Orion: ...\"\n// Every occurrence of the old sentence gets the same new wording, so we\n// search the whole body for the exact old text and replace each match.\n\nconst oldText =\n  \"Waktu Kampanye rasi bintang Orion 2022: 16-25 Januari, 14-23 Februari, 14-24 Maret\";\nconst newText =\n  \"Waktu Kampanye 2022 untuk rasi bintang Orion: 16-25 Januari, 14-23 Februari, 14-24 Maret\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The Orion activity guide repeats a campaign-dates sentence several times\n# throughout the document (once per language/section copy). The author\n# reworded it from:\n#   \"Waktu Kampanye rasi bintang Orion 2022: ...\"\n# to:\n#   \"Waktu Kampanye 2022 untuk rasi bintang Orion: ...\"\n# Every occurrence of the old sentence gets the same new wording, so we run\n# a Find/Replace All over the whole document content.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"Waktu Kampanye rasi bintang Orion 2022: 16-25 Januari, 14-23 Februari, 14-24 Maret\"\n$find.Replacement.Text = \"Waktu Kampanye 2022 untuk rasi bintang Orion: 16-25 Januari, 14-23 Februari, 14-24 Maret\"\n\n# wdFindContinue (1) wrap, wdReplaceAll (2) replace every match in the range\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
